# Update sheet MAIN: append 50 new (LDAP) user rows (rows 510-559),
# matching commit "Update Pertanggal 24 Juli 2023 10:54 WIB".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("(LDAP) wisnu.andra", "User_WisnuAndra", 4000000000509),
    @("(LDAP) khamim", "User_Khamim", 4000000000510),
    @("(LDAP) slametr", "User_Slametr", 4000000000511),
    @("(LDAP) adythiaa", "User_Adythiaa", 4000000000512),
    @("(LDAP) aguss", "User_Aguss", 4000000000513),
    @("(LDAP) rafi.artman", "User_RafiArtman", 4000000000514),
    @("(LDAP) azisp", "User_Azisp", 4000000000515),
    @("(LDAP) heryanto", "User_Heryanto", 4000000000516),
    @("(LDAP) imran", "User_Imran", 4000000000517),
    @("(LDAP) riza", "User_Riza", 4000000000518),
    @("(LDAP) ronny.anindika", "User_RonnyAnindika", 4000000000519),
    @("(LDAP) wawan.kusworo", "User_WawanKusworo", 4000000000520),
    @("(LDAP) agus.budi", "User_AgusBudi", 4000000000521),
    @("(LDAP) samta.harahap", "User_SamtaHarahap", 4000000000522),
    @("(LDAP) fuzi.mafrozi", "User_FuziMafrozi", 4000000000523),
    @("(LDAP) novizan", "User_Novizan", 4000000000524),
    @("(LDAP) yusuf.fathurahman", "User_YusufFathurahman", 4000000000525),
    @("(LDAP) muhammad.sholikhun", "User_MuhammadSholikhun", 4000000000526),
    @("(LDAP) vingky", "User_Vingky", 4000000000527),
    @("(LDAP) admin.q180", "User_AdminQ180", 4000000000528),
    @("(LDAP) restu.dwi", "User_RestuDwi", 4000000000529),
    @("(LDAP) nikkon.septian", "User_NikkonSeptian", 4000000000530),
    @("(LDAP) dian.setiawan", "User_DianSetiawan", 4000000000531),
    @("(LDAP) muhammad.syarifudin", "User_MuhammadSyarifudin", 4000000000532),
    @("(LDAP) oqi.suhaqi", "User_OqiSuhaqi", 4000000000533),
    @("(LDAP) wardah", "User_Wardah", 4000000000534),
    @("(LDAP) indra.wijaya", "User_IndraWijaya", 4000000000535),
    @("(LDAP) nadiah.rizkiah", "User_NadiahRizkiah", 4000000000536),
    @("(LDAP) istikaro.fauziah", "User_IstikaroFauziah", 4000000000537),
    @("(LDAP) bagus.isdiantara", "User_BagusIsdiantara", 4000000000538),
    @("(LDAP) cahyana", "User_Cahyana", 4000000000539),
    @("(LDAP) ahmad.yunadi", "User_AhmadYunadi", 4000000000540),
    @("(LDAP) muhammad.lukbani", "User_MuhammadLukbani", 4000000000541),
    @("(LDAP) wahyu.teluk", "User_WahyuTeluk", 4000000000542),
    @("(LDAP) denny.achmad", "User_DennyAchmad", 4000000000543),
    @("(LDAP) irma.maulidawati", "User_IrmaMaulidawati", 4000000000544),
    @("(LDAP) wulanraniasih", "User_Wulanraniasih", 4000000000545),
    @("(LDAP) procurement.admin", "User_ProcurementAdmin", 4000000000546),
    @("(LDAP) ahmad.choerul", "User_AhmadChoerul", 4000000000547),
    @("(LDAP) ferdian.kriswantoro", "User_FerdianKriswantoro", 4000000000548),
    @("(LDAP) zeinurani", "User_Zeinurani", 4000000000549),
    @("(LDAP) yogi.perbangkara", "User_YogiPerbangkara", 4000000000550),
    @("(LDAP) teluknaga", "User_Teluknaga", 4000000000551),
    @("(LDAP) rizal.amri", "User_RizalAmri", 4000000000552),
    @("(LDAP) gilang.setiawan", "User_GilangSetiawan", 4000000000553),
    @("(LDAP) fabrian.danang", "User_FabrianDanang", 4000000000554),
    @("(LDAP) asep.mulyana", "User_AsepMulyana", 4000000000555),
    @("(LDAP) dede.hartanto", "User_DedeHartanto", 4000000000556),
    @("(LDAP) yogo", "User_Yogo", 4000000000557),
    @("(LDAP) suparji", "User_Suparji", 4000000000558)
)

$startRow = 510

# --- Set the base style (wrap/left/center) on the first new cell in column B,
# then propagate it via format-only paste so the style list doesn't
# accumulate transient/orphaned cellXfs entries. ---
$firstB = $ws.Cells.Item($startRow, 2)
$firstB.Value = $data[0][0]
$firstB.WrapText = $true
$firstB.HorizontalAlignment = -4131
$firstB.VerticalAlignment = -4108

$firstB.Copy()
$lastRow = $startRow + $data.Count - 1
if ($data.Count -gt 1) {
  $ws.Range($ws.Cells.Item($startRow + 1, 2), $ws.Cells.Item($lastRow, 2)).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Fill in the LDAP name / username / PID for every new row. ---
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]

  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]

  $eCell = $ws.Cells.Item($r, 5)
  $eCell.Value = $row[2]
  $eCell.NumberFormat = "0"
  $eCell.Interior.Color = 5296274
}

# --- Column E got a bit narrower after the update. ---
$ws.Columns.Item(5).ColumnWidth = 11.3

# --- Restore the view: scrolled down with E517 selected. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 508
$win.ScrollColumn = 1
$ws.Range("E517").Select()
